$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# v0.4 restructure: columns B-D become the "common options" shared by every
# case, columns E-G hold the options specific to this case
# (target = common options, value = options by case).
# ---------------------------------------------------------------------------

# --- Row 1 (step/command names) ---
$ws.Range("C1").Value = "doubleClickAt"
$ws.Range("D1").Value = "sendKeys"
$ws.Range("E1").Value = "open"
$ws.Range("F1").Value = "doubleClickAt"
$ws.Range("G1").Value = "wait"

# --- Row 2 (targets / values) ---
$ws.Range("B2").Value = "http://127.0.0.1:9001/doubleClickAt/"
$ws.Range("C2").Value = "id=btn1"
$ws.Range("D2").Value = "xpath=//body"
$ws.Range("E2").Value = "http://127.0.0.1:9001/doubleClickAt/"
$ws.Range("F2").Value = '{"target": "id=btn1","value":"1,1"}'

# --- Row 3 ---
$ws.Range("D3").Value = '${ENTER_KEYS}'
$ws.Range("G3").Value = 2000

# ---------------------------------------------------------------------------
# Formatting: copy the already-defined visual styles onto the newly
# populated cells that reuse an existing look (do this before changing
# C2/D2's own look below, since F2/G2 borrow C2/D2's ORIGINAL style).
# ---------------------------------------------------------------------------
$xlPasteFormats = -4122

$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F1").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G1").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("B1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("D3").Copy() | Out-Null
$ws.Range("C3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("F3").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("G3").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Range("A3").Copy() | Out-Null
$ws.Range("E3").PasteSpecial($xlPasteFormats) | Out-Null

# F2/G2 keep the plain monospace "value" look the old C2/D2 cells used.
$ws.Range("C2").Copy() | Out-Null
$ws.Range("F2").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("D2").Copy() | Out-Null
$ws.Range("G2").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

# New look used only by C2 ("id=btn1"): green Sarasa Mono CL text.
$c2 = $ws.Range("C2")
$c2.Font.Name = "Sarasa Mono CL"
$c2.Font.Size = 12
$c2.Font.Color = 1539334
$c2.Borders.LineStyle = 1
$c2.VerticalAlignment = -4108

# New look used only by D2 ("xpath=//body"): green MS P Gothic text.
$d2 = $ws.Range("D2")
$d2.Font.Name = "ＭＳ Ｐゴシック"
$d2.Font.Size = 12
$d2.Font.Color = 1539334
$d2.Borders.LineStyle = 1
$d2.VerticalAlignment = -4108

# ---------------------------------------------------------------------------
# Hyperlinks: refresh B2's target and give E2 (the case-specific "open" url)
# its own hyperlink, matching the new http://127.0.0.1:9001 host.
# Updating the Address of an *existing* Hyperlink object in place (instead
# of Delete+Add) keeps B2's original style/xf untouched.
# ---------------------------------------------------------------------------
foreach ($h in $ws.Hyperlinks) {
    $h.Address = "http://127.0.0.1:9001/doubleClickAt/"
}

# E2 needs a brand new hyperlink; Hyperlinks.Add() re-stamps its own look,
# so re-apply B2's (untouched) pristine hyperlink format afterwards to land
# back on the very same style index instead of a freshly minted duplicate.
$ws.Hyperlinks.Add($ws.Range("E2"), "http://127.0.0.1:9001/doubleClickAt/") | Out-Null
$ws.Range("B2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial($xlPasteFormats) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Column widths (best-fit, like the original authoring flow) and selection.
# ---------------------------------------------------------------------------
$ws.Range("B1:F3").Columns.AutoFit() | Out-Null

$ws.Range("E10").Select() | Out-Null
